# Update scraped_at timestamps on the "snapshot" sheet (K2:K52)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("snapshot")

$ws1.Range("K2").Value = "2025-11-09T04:51:32.557554+00:00"
$ws1.Range("K3").Value = "2025-11-09T04:51:32.557619+00:00"
$ws1.Range("K4").Value = "2025-11-09T04:51:32.557648+00:00"
$ws1.Range("K5").Value = "2025-11-09T04:51:34.658854+00:00"
$ws1.Range("K6").Value = "2025-11-09T04:51:34.658869+00:00"
$ws1.Range("K7").Value = "2025-11-09T04:51:34.658877+00:00"
$ws1.Range("K8").Value = "2025-11-09T04:51:37.178423+00:00"
$ws1.Range("K9").Value = "2025-11-09T04:51:39.623657+00:00"
$ws1.Range("K10").Value = "2025-11-09T04:51:39.623688+00:00"
$ws1.Range("K11").Value = "2025-11-09T04:51:39.623708+00:00"
$ws1.Range("K12").Value = "2025-11-09T04:51:42.349622+00:00"
$ws1.Range("K13").Value = "2025-11-09T04:51:42.349651+00:00"
$ws1.Range("K14").Value = "2025-11-09T04:51:42.349670+00:00"
$ws1.Range("K15").Value = "2025-11-09T04:51:42.349688+00:00"
$ws1.Range("K16").Value = "2025-11-09T04:51:47.324466+00:00"
$ws1.Range("K17").Value = "2025-11-09T04:51:49.714643+00:00"
$ws1.Range("K18").Value = "2025-11-09T04:51:52.542628+00:00"
$ws1.Range("K19").Value = "2025-11-09T04:51:52.542657+00:00"
$ws1.Range("K20").Value = "2025-11-09T04:51:52.542676+00:00"
$ws1.Range("K21").Value = "2025-11-09T04:51:55.253319+00:00"
$ws1.Range("K22").Value = "2025-11-09T04:51:57.988689+00:00"
$ws1.Range("K23").Value = "2025-11-09T04:51:57.988720+00:00"
$ws1.Range("K24").Value = "2025-11-09T04:52:00.755599+00:00"
$ws1.Range("K25").Value = "2025-11-09T04:52:00.755629+00:00"
$ws1.Range("K26").Value = "2025-11-09T04:52:00.755648+00:00"
$ws1.Range("K27").Value = "2025-11-09T04:52:03.068364+00:00"
$ws1.Range("K28").Value = "2025-11-09T04:52:03.068407+00:00"
$ws1.Range("K29").Value = "2025-11-09T04:52:03.068428+00:00"
$ws1.Range("K30").Value = "2025-11-09T04:52:03.068438+00:00"
$ws1.Range("K31").Value = "2025-11-09T04:52:03.068445+00:00"
$ws1.Range("K32").Value = "2025-11-09T04:52:05.527655+00:00"
$ws1.Range("K33").Value = "2025-11-09T04:52:05.527691+00:00"
$ws1.Range("K34").Value = "2025-11-09T04:52:07.619053+00:00"
$ws1.Range("K35").Value = "2025-11-09T04:52:07.619071+00:00"
$ws1.Range("K36").Value = "2025-11-09T04:52:07.619079+00:00"
$ws1.Range("K37").Value = "2025-11-09T04:52:10.053713+00:00"
$ws1.Range("K38").Value = "2025-11-09T04:52:10.053746+00:00"
$ws1.Range("K39").Value = "2025-11-09T04:52:10.053766+00:00"
$ws1.Range("K40").Value = "2025-11-09T04:52:12.060906+00:00"
$ws1.Range("K41").Value = "2025-11-09T04:52:12.060938+00:00"
$ws1.Range("K42").Value = "2025-11-09T04:52:12.060957+00:00"
$ws1.Range("K43").Value = "2025-11-09T04:52:12.060975+00:00"
$ws1.Range("K44").Value = "2025-11-09T04:52:12.060991+00:00"
$ws1.Range("K45").Value = "2025-11-09T04:52:12.061006+00:00"
$ws1.Range("K46").Value = "2025-11-09T04:52:14.516175+00:00"
$ws1.Range("K47").Value = "2025-11-09T04:52:14.516193+00:00"
$ws1.Range("K48").Value = "2025-11-09T04:52:18.605382+00:00"
$ws1.Range("K49").Value = "2025-11-09T04:52:18.605413+00:00"
$ws1.Range("K50").Value = "2025-11-09T04:52:18.605431+00:00"
$ws1.Range("K51").Value = "2025-11-09T04:52:21.024079+00:00"
$ws1.Range("K52").Value = "2025-11-09T04:52:21.024100+00:00"


# Remove the two processed rows from the "new_injured" sheet, leaving just the header
$ws3 = $wb.Worksheets.Item("new_injured")
$ws3.Rows("2:3").Delete()
